# Apply the "Europe / Italy Spring 2022 / Mavic" re-label + mobile reorder
# edit described by the commit:
#   "Changed title of Europe page, re ordered people and things index page
#    for mobile"
#
# The sheet is a little HTML-snippet generator: column A holds the (rich
# text) opening half of an <img> tag, column B holds a bare filename,
# column C holds the (rich text) closing half + alt text, and column D is
# =CONCATENATE(A,B,C). Rows 2-11 used to be a scratch list of Design photos;
# they are replaced with 10 new Italy/Mavic drone filenames, the shared
# "/Design/" and alt-text fragments are re-worded, and the now-unused rows
# 12-38 (which previously held the rest of the old photo list) are cleared
# out and shrunk back down to default row height.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. New column A / C text (the two rich-text fragments around the filename)
# ---------------------------------------------------------------------
$newA = '<div class="masonryImage"> <img src="images/Europe/Italy Spring 2022/Mavic'
$newC = '" alt="Photos from Florence, Lupaia in Tuscany, and Cinque Terre in Spring of 2022." /></div>'

# ---------------------------------------------------------------------
# 2. New filenames for column B, rows 2-11
# ---------------------------------------------------------------------
$files = @(
    "DJI_0157web.jpg",
    "DJI_0160web.jpg",
    "DJI_0163web.jpg",
    "DJI_0167web.jpg",
    "DJI_0173web.jpg",
    "DJI_0176web.jpg",
    "DJI_0178web.jpg",
    "DJI_0181web.jpg",
    "DJI_0191web.jpg",
    "DJI_0195web.jpg"
)

for ($i = 0; $i -lt $files.Length; $i++) {
    $r = 2 + $i
    $ws.Range("A$r").Value = $newA
    $ws.Range("B$r").Value = $files[$i]
    $ws.Range("C$r").Value = $newC
}

# Row 2 used to be taller (ht=80) to fit the old text; restore it to the
# same height the rest of the list rows use now that the text is shorter.
$ws.Rows.Item(2).RowHeight = 64

# ---------------------------------------------------------------------
# 3. The old list continued through row 38 - clear all of that out now
#    that only 10 rows (2-11) of data remain, then let the rows collapse
#    back down to the sheet's default height.
# ---------------------------------------------------------------------
$ws.Range("A12:D38").ClearContents()
$ws.Range("A12:D38").EntireRow.AutoFit()

# ---------------------------------------------------------------------
# 4. Move the selection to D11 (last populated data row) to match the
#    saved cursor position.
# ---------------------------------------------------------------------
$ws.Range("D11").Select()
